# Handoff transform failed for both target-language sheets (zh-cn, de-de):
#  - Status (B2) flips from "Ready for handoff" to "Handoff transform failed"
#  - Latest Handoff File (C2) is cleared - no handoff file was produced, so
#    the cell (and its hyperlink) is removed entirely
#  - Latest Handoff Datetime (D2) resets to the zero date
#  - Handoff Reason (H2) flips from "Include" to "Ignored"

$wb = $excel.ActiveWorkbook

# The Overview sheet rolls up the per-language status in columns B (zh-cn)
# and C (de-de); it shares the same "status" text as the per-language
# sheets, so it flips to the new status too.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

$sheetNames = @("zh-cn", "de-de")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Remove the hyperlink attached to C2 (the "Latest Handoff File" link)
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        if ($addr -eq '$C$2') {
            $h.Delete()
        }
    }

    $ws.Range("B2").Value = "Handoff transform failed"
    $ws.Range("C2").Clear()
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"
}
